# Auto-generated edit script: update cryptos worksheet values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.745.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.239.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.85"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.239.03"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.508"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.20"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.767.65"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.768.54"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.231.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.26"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.739"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.73"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.31"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.02"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.38"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.122"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +35.77%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.96"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.04"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.31"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.19"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.54"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "510.13"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.57"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0775"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +15.45%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.08"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.64%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0423"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.131"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.79"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.301"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.912.25"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.36"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.45"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.11%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.38"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.70%  "
